$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Parameter-table numeric results: replace the plain point-estimate text
#    with "estimate (std.err)" formatted values. Every old value below is
#    unique in the document, so a single wdReplaceAll Find.Execute per pair
#    is sufficient and safe.
# ---------------------------------------------------------------------------
$replacements = @(
    @("202079171.5", "2.02079e+08 (1.68179e+06)"),
    @("7736313.0",   "7.73644e+06 (1.77182e+05)"),
    @("-34506.4",    "-3.45074e+04 (1.30100e+03)"),
    @("437.6",       "4.37642e+02 (4.69461e+00)"),
    @("0.713",       "7.13054e-01 (1.19701e-03)"),

    @("168914955.9", "1.68915e+08 (1.34721e+06)"),
    @("4317538.2",   "4.31763e+06 (1.26037e+05)"),
    @("-1402.0",     "-1.40268e+03 (9.49444e+02)"),
    @("64.43",       "6.44335e+01 (3.09119e+00)"),

    @("126265465.6", "1.26265e+08 (3.43752e+06)"),
    @("736221.9",    "7.36231e+05 (1.30763e+05)"),
    @("13332.3",     "1.33322e+04 (1.08820e+03)"),
    @("334.65",      "3.34660e+02 (8.11936e+00)"),
    @("1.698",       "1.69844e+00 (3.21488e-03)"),

    @("120413098.3", "1.20413e+08 (2.96211e+06)"),
    @("-924004.6",   "-9.23996e+05 (1.08729e+05)"),
    @("65919.7",     "6.59196e+04 (6.66640e+02)"),
    @("-317.1",      "-3.17110e+02 (2.16245e+00)"),

    @("0.0011045",   "1.10452e-03 (3.40698e-06)"),
    @("0.0022809",   "2.28095e-03 (5.00415e-06)"),
    @("35.110340",   "3.51103e+01 (1.42411e-01)"),
    @("1.1493792",   "1.14938e+00 (2.34615e-03)"),
    @("0.4411979",   "4.41198e-01 (7.91980e-03)"),

    @("0.35196",     "4.10e-01 (1.97e-01)"),
    @("-1.6111E-3",  "-1.74e-02 (2.42e-02)"),

    @("75.561",      "7.52e+01 (8.01e+00)"),
    @("-16.635",     "-1.67e+0 (9.88e-01)")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# ---------------------------------------------------------------------------
# 2) Resize the column grids of the four delta-/gamma-efficacy coefficient
#    tables to the new widths (values are twips -- Word's PreferredWidth
#    property takes points, so divide by 20).
# ---------------------------------------------------------------------------
function Set-TableColumnWidths {
    param($table, [int[]]$widthsTwips)
    for ($c = 1; $c -le $widthsTwips.Count; $c++) {
        $col = $table.Columns.Item($c)
        $col.PreferredWidthType = 3
        $col.PreferredWidth = $widthsTwips[$c - 1] / 20.0
    }
}

# Table order in the document: 9 = delta a_i, 10 = delta b_i,
# 11 = gamma a_i, 12 = gamma b_i.
Set-TableColumnWidths $d.Tables.Item(9)  @(1765, 1658, 1605, 1444, 1444)
Set-TableColumnWidths $d.Tables.Item(10) @(2160, 2029, 1963, 1767)
Set-TableColumnWidths $d.Tables.Item(11) @(1560, 1618, 1560, 1618, 1560)
Set-TableColumnWidths $d.Tables.Item(12) @(1944, 2015, 1944, 2015)
